$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update task durations (C3:C13), applying the "0" number format to each
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 30
$ws.Range("C5").Value = 30
$ws.Range("C6").Value = 20
$ws.Range("C7").Value = 90
$ws.Range("C8").Value = 120
$ws.Range("C9").Value = 120
$ws.Range("C10").Value = 45
$ws.Range("C11").Value = 120
$ws.Range("C12").Value = 45
$ws.Range("C13").Value = 120

# Apply the integer number format to the duration column C3:C13
$ws.Range("C3:C13").NumberFormat = "0"

# New blank rows for additional promotion / time-management tasks
$ws.Range("C14").NumberFormat = "0"
$ws.Range("C15").NumberFormat = "0"

# Total row keeps the same formula; recalculates automatically
$ws.Range("C16").NumberFormat = "0"

# Move the active selection to D8
$ws.Range("D8").Select()

$wb.Save()
